$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.194.70'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.927.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.89'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7131'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3203'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.36'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07089'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7899'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07941'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.926.98'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.363'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.78'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.59'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.209.05'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '255.66'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008001'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.756'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.181.31'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.808'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.513'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.60'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.01'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.270'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1264'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.355'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.527'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.383'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05148'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.266'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7431'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.761'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01948'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.799'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '77.51'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.328'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4478'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.973'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8412'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.57'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.693'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.417'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.883'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +9.08%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06107'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.56%  '
